# edit.ps1 — applies the data202402.xlsx update described by the commit diff.
#
# Summary of the edit:
#  1. Sheet "部门情况202402"   : two small corrections (F7/G7, F9/G9).
#  2. Sheet "经办人情况202402" : one small correction (L20/M20).
#  3. Sheet "个人经营贷202402" : a brand-new product row ("个人经营贷") is
#     inserted at row 2 (pushing every existing row down by one), and three
#     rows get updated totals to account for the newly reported activity
#     (瑞e惠-平安普惠, 瑞商贷, and the grand Total row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a value into a cell as plain text (matches the workbook's
# convention of storing these particular numbers as inlineStr/text rather
# than as numeric cells), without leaving behind a new cell style.
# ---------------------------------------------------------------------------
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# ---------------------------------------------------------------------------
# 1) 部门情况202402 — rows 7 & 9, columns F/G (stored as text numbers)
# ---------------------------------------------------------------------------
$wsDept = $wb.Worksheets.Item("部门情况202402")

Set-TextValue $wsDept.Range("F7") "266742.46"
Set-TextValue $wsDept.Range("G7") "2925.00"

Set-TextValue $wsDept.Range("F9") "101211.20"
Set-TextValue $wsDept.Range("G9") "6583.00"

# ---------------------------------------------------------------------------
# 2) 经办人情况202402 — row 20, columns L/M (stored as text numbers)
# ---------------------------------------------------------------------------
$wsAgent = $wb.Worksheets.Item("经办人情况202402")

Set-TextValue $wsAgent.Range("L20") "71265.10"
Set-TextValue $wsAgent.Range("M20") "1952.00"

# ---------------------------------------------------------------------------
# 3) 个人经营贷202402 — insert the new "个人经营贷" row at row 2, then patch
#    the three rows whose running totals move.
# ---------------------------------------------------------------------------
$wsLoan = $wb.Worksheets.Item("个人经营贷202402")

# Insert a blank row above the current row 2 — everything below (old rows
# 2-14) shifts down to rows 3-15, and the sheet dimension grows to A1:M15
# automatically.
$wsLoan.Rows.Item(2).Insert()

# The newly inserted row inherits the header row's bold/bordered style;
# strip that back to the plain/no-style look used by every other data row
# before filling in the new product's figures.
$wsLoan.Range("A2:M2").ClearFormats()
$wsLoan.Range("A2").Value = "个人经营贷"
$wsLoan.Range("B2").Value = 1
$wsLoan.Range("C2").Value = 60.45
$wsLoan.Range("D2").Value = 2
$wsLoan.Range("E2").Value = 245
$wsLoan.Range("F2").Value = 11
$wsLoan.Range("G2").Value = 0
$wsLoan.Range("H2").Value = 0
$wsLoan.Range("I2").Value = 0
$wsLoan.Range("J2").Value = 0
$wsLoan.Range("K2").Value = 60.45
$wsLoan.Range("L2").Value = 60.45
$wsLoan.Range("M2").Value = 60.45

# Row 7 (after the shift) = 瑞e惠-平安普惠 — new accounts/balances reported.
$wsLoan.Range("B7").Value = 1952
$wsLoan.Range("C7").Value = 71265.10000000001
$wsLoan.Range("D7").Value = 1952
$wsLoan.Range("E7").Value = 93350.89999999999
$wsLoan.Range("K7").Value = 36.51
$wsLoan.Range("L7").Value = 1475.07

# Row 10 (after the shift) = 瑞商贷 — new accounts/balances reported.
$wsLoan.Range("B10").Value = 6328
$wsLoan.Range("C10").Value = 90153.00999999999
$wsLoan.Range("D10").Value = 7748
$wsLoan.Range("E10").Value = 124861.7
$wsLoan.Range("L10").Value = 6357.11

# Row 15 (after the shift) = Total — recomputed grand totals.
$wsLoan.Range("B15").Value = 9657
$wsLoan.Range("C15").Value = 391787.6
$wsLoan.Range("D15").Value = 11499
$wsLoan.Range("E15").Value = 450722.44
$wsLoan.Range("F15").Value = 112.6
$wsLoan.Range("K15").Value = 1236.27
$wsLoan.Range("L15").Value = 10373.2
$wsLoan.Range("M15").Value = 4696.46
